$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet (reflected in workbook.xml <sheet name=.../> and app.xml titles)
$ws.Name = "Through 2022-10-31"

# Update the header label for the "2022" column (shared string used by I1)
$ws.Range("I1").Value = "2022 (through 10-31)"

# Update November carjacking count (new data point)
$ws.Range("I11").Value = 125

# Update the Total for column I accordingly
$ws.Range("I14").Value = 1401
